$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.203030303030303
$ws.Range("C2").Value = 0.5606060606060606
$ws.Range("J2").Value = 0.006060606060606061
$ws.Range("P2").Value = 0.1242424242424242
$ws.Range("S2").Value = 0.1060606060606061
$ws.Range("B3").Value = 0.005154639175257732
$ws.Range("C3").Value = 0.03608247422680412
$ws.Range("J3").Value = 0.07216494845360824
$ws.Range("P3").Value = 0.7010309278350515
$ws.Range("S3").Value = 0.1855670103092784
$ws.Range("B6").Value = 0.08502024291497975
$ws.Range("D6").Value = 0.008097165991902834
$ws.Range("F6").Value = 0.07692307692307693
$ws.Range("J6").Value = 0.2064777327935223
$ws.Range("O6").Value = 0.01619433198380567
$ws.Range("Q6").Value = 0.1781376518218623
$ws.Range("R6").Value = 0.09716599190283401
$ws.Range("S6").Value = 0.3319838056680162
$ws.Range("B7").Value = 0.1666666666666667
$ws.Range("D7").Value = 0.00641025641025641
$ws.Range("E7").Value = 0.00641025641025641
$ws.Range("F7").Value = 0.08333333333333333
$ws.Range("J7").Value = 0.108974358974359
$ws.Range("O7").Value = 0.01923076923076923
$ws.Range("Q7").Value = 0.1666666666666667
$ws.Range("R7").Value = 0.0641025641025641
$ws.Range("S7").Value = 0.3782051282051282
$ws.Range("B8").Value = 0.08616780045351474
$ws.Range("D8").Value = 0.01360544217687075
$ws.Range("F8").Value = 0.07256235827664399
$ws.Range("J8").Value = 0.1292517006802721
$ws.Range("O8").Value = 0.0272108843537415
$ws.Range("Q8").Value = 0.2063492063492063
$ws.Range("R8").Value = 0.08616780045351474
$ws.Range("S8").Value = 0.3786848072562358
$ws.Range("B9").Value = 0.09523809523809523
$ws.Range("D9").Value = 0.01058201058201058
$ws.Range("F9").Value = 0.05291005291005291
$ws.Range("J9").Value = 0.126984126984127
$ws.Range("O9").Value = 0.01058201058201058
$ws.Range("Q9").Value = 0.1851851851851852
$ws.Range("R9").Value = 0.07407407407407407
$ws.Range("S9").Value = 0.4444444444444444
$ws.Range("B10").Value = 0.1239078633836378
$ws.Range("D10").Value = 0.01906274821286736
$ws.Range("E10").Value = 0.00238284352660842
$ws.Range("F10").Value = 0.08101667990468626
$ws.Range("J10").Value = 0.1262907069102462
$ws.Range("O10").Value = 0.01826846703733121
$ws.Range("Q10").Value = 0.193010325655282
$ws.Range("R10").Value = 0.06433677521842732
$ws.Range("S10").Value = 0.3717235901509134
$ws.Range("G11").Value = 0.1367521367521368
$ws.Range("J11").Value = 0.08547008547008547
$ws.Range("K11").Value = 0.1794871794871795
$ws.Range("L11").Value = 0.5811965811965812
$ws.Range("S11").Value = 0.0170940170940171
$ws.Range("G12").Value = 0.7246376811594203
$ws.Range("J12").Value = 0.2391304347826087
$ws.Range("K12").Value = 0.01449275362318841
$ws.Range("L12").Value = 0.01449275362318841
$ws.Range("S12").Value = 0.007246376811594203
$ws.Range("G13").Value = 0.6976744186046512
$ws.Range("J13").Value = 0.2790697674418605
$ws.Range("S13").Value = 0.02325581395348837
$ws.Range("F15").Value = 0.01818181818181818
$ws.Range("H15").Value = 0.1545454545454545
$ws.Range("I15").Value = 0.1045454545454545
$ws.Range("J15").Value = 0.3227272727272728
$ws.Range("K15").Value = 0.06363636363636363
$ws.Range("M15").Value = 0.01818181818181818
$ws.Range("O15").Value = 0.06818181818181818
$ws.Range("S15").Value = 0.25
$ws.Range("F16").Value = 0.005235602094240838
$ws.Range("H16").Value = 0.1884816753926702
$ws.Range("I16").Value = 0.0418848167539267
$ws.Range("J16").Value = 0.4397905759162304
$ws.Range("K16").Value = 0.1047120418848168
$ws.Range("M16").Value = 0.01570680628272251
$ws.Range("N16").Value = 0.005235602094240838
$ws.Range("O16").Value = 0.04712041884816754
$ws.Range("S16").Value = 0.1518324607329843
$ws.Range("F17").Value = 0.02059496567505721
$ws.Range("H17").Value = 0.1739130434782609
$ws.Range("I17").Value = 0.09610983981693363
$ws.Range("J17").Value = 0.4393592677345537
$ws.Range("K17").Value = 0.07093821510297482
$ws.Range("M17").Value = 0.02517162471395881
$ws.Range("O17").Value = 0.06864988558352403
$ws.Range("S17").Value = 0.1052631578947368
$ws.Range("F18").Value = 0.02994011976047904
$ws.Range("H18").Value = 0.1616766467065868
$ws.Range("I18").Value = 0.08982035928143713
$ws.Range("J18").Value = 0.4730538922155689
$ws.Range("K18").Value = 0.09580838323353294
$ws.Range("M18").Value = 0.01796407185628742
$ws.Range("O18").Value = 0.04191616766467066
$ws.Range("S18").Value = 0.08982035928143713
$ws.Range("F19").Value = 0.01493775933609959
$ws.Range("H19").Value = 0.2232365145228216
$ws.Range("I19").Value = 0.08298755186721991
$ws.Range("J19").Value = 0.3842323651452282
$ws.Range("K19").Value = 0.08713692946058091
$ws.Range("M19").Value = 0.02157676348547718
$ws.Range("N19").Value = 0.0008298755186721991
$ws.Range("O19").Value = 0.08132780082987552
$ws.Range("S19").Value = 0.1037344398340249
